# Rebuild paragraph 1 (splitting " by Aghila" into " by " + a spell-checked
# "Aghila" run, wrapped in proofErr spellStart/spellEnd markers around the
# two "unknown" words/run-group), then append the two new paragraphs
# described in the commit ("changed docx to 2 line file"):
#   - "Hey there my name is Aghila im learning git" (with "im" flagged by
#     proofErr spellStart/spellEnd)
#   - a trailing empty paragraph
#
# The runtime merges a paragraph's runs whenever any Range.Delete/Text
# mutation touches it, which would destroy the original w:rsidR bookkeeping
# on the untouched runs. To avoid that, the fully-formed replacement
# paragraph (with the original w:rsidR/paraId attributes preserved) is
# inserted fresh via Range.InsertXML, and only then is the old paragraph 1
# removed.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

$para1Xml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="49851136" w14:textId="309569FD" w:rsidR="00EE35DA" w:rsidRPr="00A71CA6" w:rsidRDefault="00A71CA6">' +
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Hello world</w:t></w:r>' +
  '<w:r w:rsidR="00CA655F"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> by </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Aghila</w:t></w:r>' +
  '<w:r w:rsidR="003D2142"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>aaaaa</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>'

$para2Xml = '<w:p ' + $wNs + '>' +
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Hey there my name is Aghila </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>im</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> learning git</w:t></w:r>' +
  '</w:p>'

$para3Xml = '<w:p ' + $wNs + '><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'

# Insert the three replacement/new paragraphs at the very end of the body.
$insertionPoint = $d.Content
$insertionPoint.Collapse(0)
$insertionPoint.InsertXML($para1Xml + $para2Xml + $para3Xml)

# Now drop the original (stale) first paragraph; this leaves the freshly
# inserted paragraphs as paragraphs 1, 2 and 3.
$d.Paragraphs(1).Range.Delete()
